$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 348.66666
$ws.Range("K33").Value = 348.66666
$ws.Range("M33").Value = -119.66666
$ws.Range("I33").Value = 348.66666
$ws.Range("K43").Value = 3064.4285
$ws.Range("M43").Value = -2995.4285
$ws.Range("I43").Value = 3064.4285
$ws.Range("H43").Value = 3368.875
$ws.Range("K74").Value = 4131
$ws.Range("M74").Value = -3195
$ws.Range("I74").Value = 4131
$ws.Range("H74").Value = 4131
$ws.Range("K77").Value = 20655
$ws.Range("M77").Value = -15975
$ws.Range("I77").Value = 4131
$ws.Range("H77").Value = 4131
$ws.Range("H111").Value = 839.1429000000001
$ws.Range("N111").Value = -8018.9999
$ws.Range("L111").Value = 1884.9999
$ws.Range("J111").Value = 628.3333
$ws.Range("K132").Value = 14994.345
$ws.Range("M132").Value = -12464.345
$ws.Range("I132").Value = 4998.115
$ws.Range("H132").Value = 4946.185
$ws.Range("K137").Value = 6319.1535
$ws.Range("M137").Value = -3769.1535
$ws.Range("I137").Value = 2106.3845
$ws.Range("H137").Value = 2106.3845

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J2").Value = 1800
$ws.Range("L2").Value = 1800
$ws.Range("N2").Value = -2026
$ws.Range("I2").Value = 1782.3334
$ws.Range("K2").Value = 1782.3334
$ws.Range("M2").Value = -1669.3334
$ws.Range("H2").Value = 1784.8572
$ws.Range("H32").Value = 4265.1665
$ws.Range("N32").Value = -15207
$ws.Range("L32").Value = 14633
$ws.Range("J32").Value = 14633
$ws.Range("K32").Value = 3113.1853
$ws.Range("M32").Value = -2826.1853
$ws.Range("I32").Value = 3113.1853
$ws.Range("I34").Value = 5000
$ws.Range("H34").Value = 5000
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = -4729
$ws.Range("I61").Value = 2232.75
$ws.Range("K61").Value = 2232.75
$ws.Range("M61").Value = -2020.75
$ws.Range("H61").Value = 2232.75
$ws.Range("K110").Value = 13106
$ws.Range("I110").Value = 13106
$ws.Range("M110").Value = -11061
$ws.Range("H110").Value = 10624.6
$ws.Range("L116").Value = 1800
$ws.Range("J116").Value = 1800
$ws.Range("K116").Value = 1782.3334
$ws.Range("M116").Value = 511.6666
$ws.Range("I116").Value = 1782.3334
$ws.Range("H116").Value = 1784.8572
$ws.Range("N116").Value = -6388
$ws.Range("J119").Value = 15218
$ws.Range("L119").Value = 15218
$ws.Range("N119").Value = -24894
$ws.Range("H119").Value = 15218
$ws.Range("K136").Value = 6698.25
$ws.Range("M136").Value = -4148.25
$ws.Range("I136").Value = 2232.75
$ws.Range("H136").Value = 2232.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M3").Value = -1668.3334
$ws.Range("N3").Value = -2028
$ws.Range("H3").Value = 1784.8572
$ws.Range("J3").Value = 1800
$ws.Range("L3").Value = 1800
$ws.Range("I3").Value = 1782.3334
$ws.Range("K3").Value = 1782.3334
$ws.Range("J20").Value = 3863.6
$ws.Range("L20").Value = 3863.6
$ws.Range("I20").Value = 1337
$ws.Range("K20").Value = 1337
$ws.Range("H20").Value = 2126.5625
$ws.Range("M20").Value = -1090
$ws.Range("N20").Value = -4357.6
$ws.Range("J86").Value = 360
$ws.Range("K86").Value = 1607.6
$ws.Range("M86").Value = -484.5999999999999
$ws.Range("I86").Value = 1607.6
$ws.Range("H86").Value = 1548.1904
$ws.Range("N86").Value = -2606
$ws.Range("L86").Value = 360
$ws.Range("L89").Value = 1800
$ws.Range("J89").Value = 360
$ws.Range("K89").Value = 8038
$ws.Range("M89").Value = -2422
$ws.Range("I89").Value = 1607.6
$ws.Range("H89").Value = 1548.1904
$ws.Range("N89").Value = -13032

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K36").Value = 6000
$ws.Range("M36").Value = -5612
$ws.Range("I36").Value = 6000
$ws.Range("H36").Value = 6000
$ws.Range("K40").Value = 6000
$ws.Range("I40").Value = 6000
$ws.Range("M40").Value = -5840
$ws.Range("H40").Value = 6000
$ws.Range("H94").Value = 1178.6
$ws.Range("N94").Value = -1583
$ws.Range("J94").Value = 681
$ws.Range("L94").Value = 681
$ws.Range("J107").Value = 1369.3334
$ws.Range("L107").Value = 1369.3334
$ws.Range("H107").Value = 1130.3846
$ws.Range("N107").Value = -5209.3334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K51").Value = 5998.5
$ws.Range("M51").Value = -5538.5
$ws.Range("I51").Value = 1999.5
$ws.Range("H51").Value = 1999.5
$ws.Range("H113").Value = 2029.1333
$ws.Range("N113").Value = -10449.5
$ws.Range("L113").Value = 6109.5
$ws.Range("J113").Value = 2036.5
$ws.Range("K113").Value = 5998.9998
$ws.Range("M113").Value = -3828.9998
$ws.Range("I113").Value = 1999.6666
$ws.Range("K132").Value = 9102.857399999999
$ws.Range("M132").Value = -6572.857399999999
$ws.Range("I132").Value = 1011.4286
$ws.Range("H132").Value = 1122.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I97").Value = 394.15
$ws.Range("K97").Value = 394.15
$ws.Range("M97").Value = 101.85
$ws.Range("H97").Value = 394.15
$ws.Range("H102").Value = 1492.3846
$ws.Range("N102").Value = -6244
$ws.Range("J102").Value = 3000
$ws.Range("L102").Value = 3000

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6539.5454
$ws.Range("K7").Value = 6539.5454
$ws.Range("I7").Value = 6539.5454
$ws.Range("M7").Value = -6427.5454
$ws.Range("I22").Value = 1444
$ws.Range("M22").Value = -1149
$ws.Range("N22").Value = -2272.5
$ws.Range("H22").Value = 1563.25
$ws.Range("J22").Value = 1682.5
$ws.Range("K22").Value = 1444
$ws.Range("L22").Value = 1682.5
$ws.Range("H27").Value = 1563.25
$ws.Range("N27").Value = -1896.5
$ws.Range("J27").Value = 1682.5
$ws.Range("L27").Value = 1682.5
$ws.Range("I27").Value = 1444
$ws.Range("K27").Value = 1444
$ws.Range("M27").Value = -1337
$ws.Range("K35").Value = 1200
$ws.Range("M35").Value = -864
$ws.Range("I35").Value = 1200
$ws.Range("H35").Value = 1200
$ws.Range("K40").Value = 1941
$ws.Range("I40").Value = 1941
$ws.Range("M40").Value = -1805
$ws.Range("H40").Value = 1941
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("M126").Value = -17148.6362
$ws.Range("H126").Value = 6539.5454
$ws.Range("K126").Value = 19618.6362
$ws.Range("I126").Value = 6539.5454
$ws.Range("H128").Value = 79283.14
$ws.Range("N128").Value = -89243.14
$ws.Range("J128").Value = 79283.14
$ws.Range("L128").Value = 79283.14
$ws.Range("N130").Value = -76705.836
$ws.Range("H130").Value = 66665.836
$ws.Range("L130").Value = 66665.836
$ws.Range("J130").Value = 66665.836
$ws.Range("K132").Value = 12563.334
$ws.Range("M132").Value = -10033.334
$ws.Range("I132").Value = 4187.778
$ws.Range("H132").Value = 4116.1665

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K107").Value = 765.75
$ws.Range("I107").Value = 255.25
$ws.Range("M107").Value = 1154.25
$ws.Range("H107").Value = 2203.8
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
$ws.Range("H119").Value = 50000
$ws.Range("M126").Value = -9487.454000000002
$ws.Range("H126").Value = 3966.8572
$ws.Range("K126").Value = 11957.454
$ws.Range("I126").Value = 3985.818
$ws.Range("K136").Value = 12853.3842
$ws.Range("M136").Value = -10303.3842
$ws.Range("I136").Value = 4284.4614
$ws.Range("H136").Value = 7093.0347
$ws.Range("K141").Value = 50000
$ws.Range("M141").Value = -44820
$ws.Range("I141").Value = 50000
$ws.Range("H141").Value = 50000
